# "add info in excel files"
# Fill in the remaining "types" problem rows (6-17), fix the duration values
# for the first two data rows, and center the duration column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("types")

$apos = [char]0x2019

# --- fix existing duration values (row 2 and row 3) ---
$ws.Range("C2").Value = 10
$ws.Range("C3").Value = 20

# --- seed the shared-string table with the new "description" text in the
#     same order it was originally authored (not plain row order) ---
$ws.Range("B6").Value  = "Switches of light not working"
$ws.Range("B7").Value  = "Turned off randomly"
$ws.Range("B8").Value  = "Remote doesn" + $apos + "t work properly"
$ws.Range("B12").Value = "Buttons not working"
$ws.Range("B15").Value = "Turns on by itself"
$ws.Range("B9").Value  = "Not heating"
$ws.Range("B10").Value = "Not cooling"
$ws.Range("B13").Value = "Not working"
$ws.Range("B14").Value = "Leaking"
$ws.Range("B11").Value = "Weird Smell"
$ws.Range("B16").Value = "No lights"

# --- "type" column for the new rows ---
$ws.Range("A6").Value  = "regular"
$ws.Range("A7").Value  = "critical"
$ws.Range("A8").Value  = "critical"
$ws.Range("A9").Value  = "regular"
$ws.Range("A10").Value = "regular"
$ws.Range("A11").Value = "regular"
$ws.Range("A12").Value = "critical"
$ws.Range("A13").Value = "critical"
$ws.Range("A14").Value = "critical"
$ws.Range("A15").Value = "critical"
$ws.Range("A16").Value = "regular"

# --- "duration" column for the new rows (row 17 only has a duration) ---
$ws.Range("C6").Value  = 10
$ws.Range("C7").Value  = 50
$ws.Range("C8").Value  = 20
$ws.Range("C9").Value  = 20
$ws.Range("C10").Value = 20
$ws.Range("C11").Value = 20
$ws.Range("C12").Value = 40
$ws.Range("C13").Value = 50
$ws.Range("C14").Value = 40
$ws.Range("C15").Value = 30
$ws.Range("C16").Value = 30
$ws.Range("C17").Value = 20

# --- make sure every cell in the duration column (header included) carries
#     the same formatting the original data cells had, then center it ---
$ws.Range("C2").Copy()
$ws.Range("C1:C17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("C1:C17").HorizontalAlignment = -4108

$ws.Range("C18").Select()
